# SIE_mockups.pptx edit: reorder "Desenvolvedores" slide earlier, and turn the
# "Alterar e-mail" mockup into a combined "old email / new email" mockup by
# relabeling the existing fields and adding the "E-mail novo" input row.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide reorder: move the "Desenvolvedores" slide (SlideID 258) from its
#    current position (8) to right after slide 3 / before the old slide4
#    (position 4), matching the new <p:sldIdLst> order:
#    256, 269, 257, 258, 260, 261, 263, 262, 259, 266, 265, 267, 264, 268
# ---------------------------------------------------------------------------
$moved = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    if ($p.Slides.Item($i).SlideID -eq 258) {
        $moved = $p.Slides.Item($i)
        break
    }
}
$moved.MoveTo(4)

# ---------------------------------------------------------------------------
# 2) Find the "Alterar e-mail" slide (SlideID 263) and edit its contents.
# ---------------------------------------------------------------------------
$emailSlide = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    if ($p.Slides.Item($i).SlideID -eq 263) {
        $emailSlide = $p.Slides.Item($i)
        break
    }
}

# Relabel the existing "new e-mail" prompt/box -> "old e-mail" prompt/box.
for ($i = 1; $i -le $emailSlide.Shapes.Count; $i++) {
    $shp = $emailSlide.Shapes.Item($i)
    if ($shp.HasTextFrame -eq -1) {
        $t = $shp.TextFrame.TextRange.Text
        if ($t -eq "Novo e-mail:") {
            $shp.TextFrame.TextRange.Text = "E-mail antigo:"
        } elseif ($t -eq "(inserir novo e-mail)") {
            $shp.TextFrame.TextRange.Text = "(inserir e-mail antigo)"
        }
    }
}

# Add the new "E-mail novo:" label textbox (plain, no fill/border).
$lbl = $emailSlide.Shapes.AddTextbox(1, 2345328 / 12700, 3849713 / 12700, 5870750 / 12700, 369332 / 12700)
$lbl.TextFrame.TextRange.Text = "E-mail novo:"

# Add the new "(inserir e-mail novo)" fill-in box (styled like the other
# input boxes on this slide: filled + bordered, centered text).
$box = $emailSlide.Shapes.AddTextbox(1, 5518814 / 12700, 3844000 / 12700, 2882721 / 12700, 369332 / 12700)
$box.TextFrame.TextRange.Text = "(inserir e-mail novo)"
$box.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$box.Fill.Solid()
$box.Fill.ForeColor.RGB = 16777215
$box.Line.Visible = $true
$box.Line.ForeColor.RGB = 0
